$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "61.463.15"
$ws.Cells.Item(2, 5).Value = "  -1.67%  "
$ws.Cells.Item(3, 4).Value = "3.355.66"
$ws.Cells.Item(3, 5).Value = "  -0.62%  "
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "400.90"
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -3.57%  "
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = "126.10"
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +7.55%  "
$ws.Cells.Item(7, 5).Value = "  +2.21%  "
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +0.03%  "
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.659"
$cell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +4.50%  "
$ws.Cells.Item(10, 5).Value = "  +1.81%  "
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = "40.94"
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +1.93%  "
$ws.Cells.Item(12, 5).Value = "  -0.99%  "
$ws.Cells.Item(13, 4).Value = "3.890.43"
$ws.Cells.Item(13, 5).Value = "  -0.40%  "
$ws.Cells.Item(14, 5).Value = "  -1.04%  "
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = "19.32"
$cell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -0.41%  "
$ws.Cells.Item(16, 4).Value = "3.352.73"
$ws.Cells.Item(16, 5).Value = "  -1.70%  "
$ws.Cells.Item(17, 4).Value = "61.434.81"
$ws.Cells.Item(17, 5).Value = "  -1.38%  "
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = "11.23"
$cell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +2.54%  "
$ws.Cells.Item(19, 5).Value = "  -0.70%  "
$ws.Cells.Item(20, 5).Value = "  +7.40%  "
$ws.Cells.Item(21, 5).Value = "  -4.34%  "
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = "80.10"
$cell.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +7.04%  "
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = "12.63"
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -0.01%  "
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = "298.88"
$cell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +0.53%  "
$ws.Cells.Item(25, 5).Value = "  -1.34%  "
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.75"
$cell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +11.19%  "
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = "8.18"
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +6.73%  "
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = "28.96"
$cell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -2.45%  "
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.49"
$cell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -6.41%  "
$ws.Cells.Item(30, 5).Value = "  -2.49%  "
$ws.Cells.Item(31, 5).Value = "  +1.05%  "
$ws.Cells.Item(32, 5).Value = "  +0.03%  "
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = "11.32"
$cell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -1.05%  "
$ws.Cells.Item(34, 5).Value = "  -1.72%  "
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = "41.02"
$cell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -4.97%  "
$ws.Cells.Item(36, 5).Value = "  -2.12%  "
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = "51.90"
$cell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -0.72%  "
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.998"
$cell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -0.02%  "
$ws.Cells.Item(39, 5).Value = "  -1.50%  "
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.91"
$cell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -6.45%  "
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = "137.24"
$cell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +2.62%  "
$ws.Cells.Item(42, 5).Value = "  +2.44%  "
$ws.Cells.Item(43, 5).Value = "  +1.08%  "
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.281"
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -2.72%  "
$ws.Cells.Item(45, 5).Value = "  +0.25%  "
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = "16.58"
$cell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +0.47%  "
$ws.Cells.Item(47, 5).Value = "  -0.45%  "
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = "21.11"
$cell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -0.85%  "
$ws.Cells.Item(49, 4).Value = "3.687.81"
$ws.Cells.Item(49, 5).Value = "  -0.50%  "
$ws.Cells.Item(50, 4).Value = "2.097.25"
$ws.Cells.Item(50, 5).Value = "  -3.50%  "
$ws.Cells.Item(51, 5).Value = "  -4.74%  "
